$p = $ppt.ActivePresentation

# --- Slide 1 / "TextBox 1": fix the "Mariko Tagawa (email), JICA volunteer" line ---
$s1 = $p.Slides.Item(1)
$creditBox = $s1.Shapes.Item(186)
$tr1 = $creditBox.TextFrame.TextRange

# Run 2 is the "Mariko Tagawa" run that carries the mailto hyperlink. Add a
# trailing space to its text and re-assert the hyperlink address; rewriting
# the Hyperlink through ActionSettings drops the legacy
# ahyp:hlinkClr extLst that used to be nested inside <a:hlinkClick>.
$nameRun = $tr1.Runs(2)
$nameRun.Text = "Mariko Tagawa "
$nameHyperlink = $nameRun.ActionSettings(1).Hyperlink
$nameHyperlink.Address = "mailto:mesa0121mesa@gmail.com"

# The remaining text " (marikotagawa@gmail.com), JICA volunteer" is spread
# across several runs (the parenthetical email plus the trailing label).
# Re-fetch the run that now follows the name (its Start shifted because of
# the space we just inserted) and collapse that whole tail into one run by
# setting the text of a Characters() range spanning it.
$tailRun = $tr1.Runs(3)
$tailLength = $tailRun.Text.Length
$tail = $tr1.Characters($tailRun.Start, $tailLength)
$tail.Text = ", JICA volunteer"

# --- Slide 12 / quiz link shape: update the displayed Google Forms URL ---
$s12 = $p.Slides.Item(12)
$linkShape = $s12.Shapes.Item(2)
$tr12 = $linkShape.TextFrame.TextRange
$urlRun = $tr12.Runs(1)
$urlRun.Text = "https://forms.gle/niLWWqzLhdUc7yoy8"
